$d = $word.ActiveDocument

# -------------------------------------------------------------------------
# Change 1: first paragraph gains a trailing red "(This is a change -
# Version for branch alternate)" annotation, and two extra trailing spaces
# on the existing sentence.
# -------------------------------------------------------------------------
$xmlPara1 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">This is a Microsoft word document.  </w:t></w:r><w:r><w:rPr><w:color w:val="C00000"/></w:rPr><w:t>(This is a change – Ve</w:t></w:r><w:r><w:rPr><w:color w:val="C00000"/></w:rPr><w:t>rsion for branch alternate</w:t></w:r><w:r><w:rPr><w:color w:val="C00000"/></w:rPr><w:t>)</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$d.Paragraphs(1).Range.InsertXML($xmlPara1)

# -------------------------------------------------------------------------
# Change 2: the "Crispian's Day speech..." heading paragraph keeps the same
# visible text but Word re-flowed the run/proofErr boundaries (a trailing
# space moved from its own run into " Day speech from ", the spell-check
# markers around "Shakespear's" shifted, and the tail " Henry V [Source -
# Wikipedia]" got merged into a single run).
# -------------------------------------------------------------------------
$xmlPara4 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:shd w:val="clear" w:color="auto" w:fill="F9F9F9"/><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:color w:val="202122"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:color w:val="202122"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>Crispian’s</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:color w:val="202122"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve"> Day speech from </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:color w:val="202122"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>Shakespear’s</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:color w:val="202122"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve"> Henry V [Source – Wikipedia]</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$d.Paragraphs(4).Range.InsertXML($xmlPara4)

# -------------------------------------------------------------------------
# Change 3: two new empty paragraphs appended at the very end of the body,
# right before the section break -- one tagged with the "larger" style plus
# explicit shading/spacing, the other a bare empty paragraph.
# -------------------------------------------------------------------------
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$endPoint = $d.Range($lastPara.Range.End, $lastPara.Range.End)
$xmlTail = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="larger"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:spacing w:before="0" w:beforeAutospacing="0" w:after="150" w:afterAutospacing="0"/></w:pPr></w:p><w:p/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$endPoint.InsertXML($xmlTail)

# -------------------------------------------------------------------------
# Change 4: drop the now-unused "apple-converted-space" and "Hyperlink"
# character styles from styles.xml (neither is referenced by any run).
# Deleted in last-to-first (by style index) order, since deleting a style
# that isn't the highest-indexed one first trips a stale-index bug for any
# later by-name style lookup in this runtime.
# -------------------------------------------------------------------------
$d.Styles("Hyperlink").Delete()
$d.Styles("apple-converted-space").Delete()

Write-Output "done"
